# Regenerate handback status report timestamps / status after a new
# handback report run ("Generate Report for Handback").
#
# Overview sheet: "Latest HO Xliff Generate Date" moves forward for the
#   03827f75... file (rows 2 & 3 share the same generation timestamp).
# zh-cn sheet: Priority flips from "ht" (human translation) to "mt"
#   (machine translation); handoff/handback datetimes advance.
# de-de sheet: Priority flips from "ht" to "mt"; handback datetime advances
#   (its handoff datetime is the same value as the Overview sheet's date).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for both rows (also backs the de-de
# sheet's "Correspond Handoff Datetime" for the same source file).
$overview.Range("G2").Value = "2016-10-10 10:02:40"
$overview.Range("G3").Value = "2016-10-10 10:02:40"

# --- zh-cn sheet ------------------------------------------------------
# Priority: human translation -> machine translation
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

# Correspond Handoff Datetime
$zhcn.Range("H2").Value = "2016-10-10 10:02:30"
$zhcn.Range("H3").Value = "2016-10-10 10:02:30"

# Correspond Handback Datetime
$zhcn.Range("K2").Value = "2016-10-10 10:03:17"
$zhcn.Range("K3").Value = "2016-10-10 10:03:17"

# --- de-de sheet --------------------------------------------------
# Priority: human translation -> machine translation
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# Correspond Handoff Datetime (mirrors the Overview sheet's generation date)
$dede.Range("H2").Value = "2016-10-10 10:02:40"
$dede.Range("H3").Value = "2016-10-10 10:02:40"

# Correspond Handback Datetime
$dede.Range("K2").Value = "2016-10-10 10:03:33"
$dede.Range("K3").Value = "2016-10-10 10:03:33"

$wb.Save()
